$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 16-18 content, plus replacement of A15 ---
# Order matters: it determines the order new strings are appended to
# sharedStrings.xml, which must match the target workbook's layout.

$ws.Range("C16").Value = "Robertson, D., J. Magnuson, S. Carpenter, and E. Stanley. 2014. North Temperate Lakes LTER Morphometry and Hypsometry data for core study lakes ver 1. Environmental Data Initiative. https://doi.org/10.6073/pasta/1d15f38aaf14110714add6230ef78bd8. Accessed 2020-06-23."
$ws.Range("B16").Value = "Hypsometry data for ME and MO"
$ws.Range("A16").Value = "lake_hypsometry"

$ws.Range("A15").Value = "level_data.rds"

$ws.Range("C17").Value = "N. Lead PI, J. Magnuson, S. Carpenter, and E. Stanley. 2020. North Temperate Lakes LTER: Snow and Ice Depth 1982 - current ver 31. Environmental Data Initiative. https://doi.org/10.6073/pasta/f4e281545ec5c5c18dc996cf652f5f8c. Accessed 2020-06-23."
$ws.Range("B17").Value = "ice thickness ME and MO"
$ws.Range("A17").Value = "winter"

$ws.Range("C18").Value = "N. Lead PI, N. LTER, J. Magnuson, S. Carpenter, and E. Stanley. 2020. North Temperate Lakes LTER: Ice Duration - Madison Lakes Area 1853 - current ver 34. Environmental Data Initiative. https://doi.org/10.6073/pasta/22a5b5f8bce193353e559918b0024f9d. Accessed 2020-06-23."
$ws.Range("B18").Value = "ice on off dates ME and MO"
$ws.Range("A18").Value = "icedates"

# --- Row heights for new rows ---
$ws.Rows.Item(16).RowHeight = 72
$ws.Rows.Item(17).RowHeight = 72
$ws.Rows.Item(18).RowHeight = 72

# --- Formatting for the citation column (C) on the new rows: wrap text,
# no explicit vertical alignment (defaults to bottom) ---
$ws.Range("C16:C18").VerticalAlignment = -4107
$ws.Range("C16:C18").WrapText = $true
$ws.Range("C16:C18").Interior.Pattern = -4142

# --- Update selection to reflect the new active cell ---
$ws.Range("A18").Select()
